$d = $word.ActiveDocument

$xmlFrag = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">Phường, mã, </w:t></w:r><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>huyện, tỉnh</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target = $d.Content
$target.Collapse(0)
$target.InsertXML($xmlFrag)
